$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column I (10/03/2023) and column J (Delta_Offerto) values for the
# rows whose "Cleared" / related figures changed (prelude of chart for MOB p.2)

$ws.Range("I2").Value = 2342.6
$ws.Range("J2").Value = -47.67941603346708

$ws.Range("I3").Value = 317
$ws.Range("J3").Value = -56.46687697160884

$ws.Range("I5").Value = 1680
$ws.Range("J5").Value = -29.22619047619047

$ws.Range("I6").Value = 1224
$ws.Range("J6").Value = -34.15032679738562

$ws.Range("I7").Value = 356
$ws.Range("J7").Value = -16.01123595505618

$ws.Range("I9").Value = 1454
$ws.Range("J9").Value = -63.61760660247593

$ws.Range("I10").Value = 208
$ws.Range("J10").Value = -54.32692307692308

$ws.Range("I12").Value = 149
$ws.Range("J12").Value = -41.61073825503355
